# Auto-generated script to apply Asura_Profits price/profit data refresh
# from the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3366.147
$ws.Range("I64").Value = 3060.6428
$ws.Range("J64").Value = 3580
$ws.Range("K64").Value = 3060.6428
$ws.Range("L64").Value = 3580
$ws.Range("M64").Value = -2812.6428
$ws.Range("N64").Value = -4076
$ws.Range("H67").Value = 3366.147
$ws.Range("I67").Value = 3060.6428
$ws.Range("J67").Value = 3580
$ws.Range("K67").Value = 3060.6428
$ws.Range("L67").Value = 3580
$ws.Range("M67").Value = -2202.6428
$ws.Range("N67").Value = -5296
$ws.Range("H70").Value = 78121.53999999999
$ws.Range("I70").Value = 501150
$ws.Range("J70").Value = 1207.2727
$ws.Range("K70").Value = 1503450
$ws.Range("L70").Value = 3621.8181
$ws.Range("M70").Value = -1503180
$ws.Range("N70").Value = -4161.8181
$ws.Range("H73").Value = 78121.53999999999
$ws.Range("I73").Value = 501150
$ws.Range("J73").Value = 1207.2727
$ws.Range("K73").Value = 1503450
$ws.Range("L73").Value = 3621.8181
$ws.Range("M73").Value = -1502514
$ws.Range("N73").Value = -5493.8181
$ws.Range("H74").Value = 4110.4546
$ws.Range("I74").Value = 3938.75
$ws.Range("J74").Value = 4208.5713
$ws.Range("K74").Value = 3938.75
$ws.Range("L74").Value = 4208.5713
$ws.Range("M74").Value = -3002.75
$ws.Range("N74").Value = -6080.5713
$ws.Range("H77").Value = 4110.4546
$ws.Range("I77").Value = 3938.75
$ws.Range("J77").Value = 4208.5713
$ws.Range("K77").Value = 19693.75
$ws.Range("L77").Value = 21042.8565
$ws.Range("M77").Value = -15013.75
$ws.Range("N77").Value = -30402.8565
$ws.Range("H137").Value = 1930.4667
$ws.Range("I137").Value = 1923.1177
$ws.Range("J137").Value = 1940.0769
$ws.Range("K137").Value = 5769.3531
$ws.Range("L137").Value = 5820.2307
$ws.Range("M137").Value = -3219.3531
$ws.Range("N137").Value = -10920.2307
$ws.Range("H138").Value = 3336535.5
$ws.Range("I138").Value = 7693596
$ws.Range("J138").Value = 4665.9116
$ws.Range("K138").Value = 23080788
$ws.Range("L138").Value = 13997.7348
$ws.Range("M138").Value = -23075648
$ws.Range("N138").Value = -24277.7348

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1278.091
$ws.Range("I74").Value = 1099.2858
$ws.Range("J74").Value = 1591
$ws.Range("K74").Value = 1099.2858
$ws.Range("L74").Value = 1591
$ws.Range("M74").Value = -225.2858000000001
$ws.Range("N74").Value = -3339
$ws.Range("H77").Value = 1278.091
$ws.Range("I77").Value = 1099.2858
$ws.Range("J77").Value = 1591
$ws.Range("K77").Value = 5496.429
$ws.Range("L77").Value = 7955
$ws.Range("M77").Value = -1128.429
$ws.Range("N77").Value = -16691
$ws.Range("H88").Value = 2982.8
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2982.8
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2982.8
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3794.8
$ws.Range("H91").Value = 2982.8
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2982.8
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2982.8
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5790.8
$ws.Range("H110").Value = 987.5833
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H115").Value = 9000
$ws.Range("I115").Value = 9000
$ws.Range("K115").Value = 9000
$ws.Range("M115").Value = -7433
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 200983.8
$ws.Range("I86").Value = 1228
$ws.Range("K86").Value = 1228
$ws.Range("M86").Value = -105
$ws.Range("H89").Value = 200983.8
$ws.Range("I89").Value = 1228
$ws.Range("K89").Value = 6140
$ws.Range("M89").Value = -524
$ws.Range("H134").Value = 2477.889
$ws.Range("I134").Value = 2334.3333
$ws.Range("K134").Value = 7002.999899999999
$ws.Range("M134").Value = -4467.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 904322.9399999999
$ws.Range("I58").Value = 1235598.6
$ws.Range("J58").Value = 844
$ws.Range("K58").Value = 1235598.6
$ws.Range("L58").Value = 844
$ws.Range("M58").Value = -1235395.6
$ws.Range("N58").Value = -1250
$ws.Range("H94").Value = 2003.7368
$ws.Range("I94").Value = 979.8
$ws.Range("J94").Value = 2369.4285
$ws.Range("K94").Value = 979.8
$ws.Range("L94").Value = 2369.4285
$ws.Range("M94").Value = -528.8
$ws.Range("N94").Value = -3271.4285
$ws.Range("H134").Value = 1198.2858
$ws.Range("I134").Value = 1084.2162
$ws.Range("J134").Value = 2042.4
$ws.Range("K134").Value = 3252.6486
$ws.Range("L134").Value = 6127.200000000001
$ws.Range("M134").Value = -717.6486000000004
$ws.Range("N134").Value = -11197.2
$ws.Range("H136").Value = 904322.9399999999
$ws.Range("I136").Value = 1235598.6
$ws.Range("J136").Value = 844
$ws.Range("K136").Value = 3706795.8
$ws.Range("L136").Value = 2532
$ws.Range("M136").Value = -3704245.8
$ws.Range("N136").Value = -7632
$ws.Range("H141").Value = 34242.8
$ws.Range("J141").Value = 32803.5
$ws.Range("L141").Value = 32803.5
$ws.Range("N141").Value = -43163.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 715736.3
$ws.Range("J12").Value = 920211
$ws.Range("L12").Value = 2760633
$ws.Range("N12").Value = -2760979
$ws.Range("H133").Value = 5288
$ws.Range("J133").Value = 6545.7144
$ws.Range("L133").Value = 19637.1432
$ws.Range("N133").Value = -29757.1432
$ws.Range("H140").Value = 2553.628
$ws.Range("I140").Value = 906.36365
$ws.Range("J140").Value = 4279.3335
$ws.Range("K140").Value = 2719.09095
$ws.Range("L140").Value = 12838.0005
$ws.Range("M140").Value = 2460.90905
$ws.Range("N140").Value = -23198.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H40").Value = 27900
$ws.Range("J40").Value = 5800
$ws.Range("L40").Value = 5800
$ws.Range("N40").Value = -6102
$ws.Range("H132").Value = 2099.6206
$ws.Range("I132").Value = 1328.9412
$ws.Range("K132").Value = 3986.8236
$ws.Range("M132").Value = -1456.8236

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 14043.4375
$ws.Range("I61").Value = 21021.1
$ws.Range("J61").Value = 2414
$ws.Range("K61").Value = 21021.1
$ws.Range("L61").Value = 2414
$ws.Range("M61").Value = -20819.1
$ws.Range("N61").Value = -2818
$ws.Range("H68").Value = 1933.3334
$ws.Range("I68").Value = 1566.6666
$ws.Range("J68").Value = 2666.6667
$ws.Range("K68").Value = 1566.6666
$ws.Range("L68").Value = 2666.6667
$ws.Range("M68").Value = -817.6666
$ws.Range("N68").Value = -4164.6667
$ws.Range("H71").Value = 1933.3334
$ws.Range("I71").Value = 1566.6666
$ws.Range("J71").Value = 2666.6667
$ws.Range("K71").Value = 7833.333000000001
$ws.Range("L71").Value = 13333.3335
$ws.Range("M71").Value = -4089.333000000001
$ws.Range("N71").Value = -20821.3335
$ws.Range("H82").Value = 1745.3704
$ws.Range("I82").Value = 1691.8823
$ws.Range("J82").Value = 1836.3
$ws.Range("K82").Value = 1691.8823
$ws.Range("L82").Value = 1836.3
$ws.Range("M82").Value = -1330.8823
$ws.Range("N82").Value = -2558.3
$ws.Range("H85").Value = 1745.3704
$ws.Range("I85").Value = 1691.8823
$ws.Range("J85").Value = 1836.3
$ws.Range("K85").Value = 1691.8823
$ws.Range("L85").Value = 1836.3
$ws.Range("M85").Value = -443.8823
$ws.Range("N85").Value = -4332.3
$ws.Range("H113").Value = 14043.4375
$ws.Range("I113").Value = 21021.1
$ws.Range("J113").Value = 2414
$ws.Range("K113").Value = 21021.1
$ws.Range("L113").Value = 2414
$ws.Range("M113").Value = -18851.1
$ws.Range("N113").Value = -6754
$ws.Range("H132").Value = 7597.0835
$ws.Range("I132").Value = 7920.875
$ws.Range("J132").Value = 6949.5
$ws.Range("K132").Value = 23762.625
$ws.Range("L132").Value = 20848.5
$ws.Range("M132").Value = -21232.625
$ws.Range("N132").Value = -25908.5
$ws.Range("H136").Value = 2252.5293
$ws.Range("I136").Value = 1932.5416
$ws.Range("J136").Value = 3020.5
$ws.Range("K136").Value = 5797.6248
$ws.Range("L136").Value = 9061.5
$ws.Range("M136").Value = -3247.6248
$ws.Range("N136").Value = -14161.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 83396
$ws.Range("J58").Value = 100094
$ws.Range("L58").Value = 100094
$ws.Range("N58").Value = -100710
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -31240
$ws.Range("H122").Value = 1965.6774
$ws.Range("I122").Value = 2024.2106
$ws.Range("J122").Value = 1873
$ws.Range("K122").Value = 6072.6318
$ws.Range("L122").Value = 5619
$ws.Range("M122").Value = -3622.6318
$ws.Range("N122").Value = -10519
